$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Prevent Excel from auto-converting the "dd-MMM-yyyy" look-alike text into
# a real date serial number; keep it as literal text like the source file.
$ws.Range("I3:I29").NumberFormat = "@"

for ($row = 3; $row -le 29; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H: PERIOD TO EXPIRE
    $hCell.Value = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # Column I: LAST UPDATE
    $iCell.Value = "04-Nov-2025"
}
